# daily auto push: 2026-02-28 18:48 UTC
# Insert two new rows of data (2026/02/28 22:00 slot and 2026/03/01 01:00 slot)
# just before the "2026/12/29" block, shifting all subsequent rows down by 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 895:896 - pushes old rows 895.. down to 897..
$ws.Rows("895:896").Insert()

# Row 895: 2026/02/28, 土, 22, 201
$ws.Range("A895").NumberFormat = "@"
$ws.Range("A895").Value = "2026/02/28"
$ws.Range("A895").Style = "Normal"
$ws.Range("B895").Value = "土"
$ws.Range("C895").Value = 22
$ws.Range("D895").Value = 201

# Row 896: 2026/03/01, 日, 1, 201
$ws.Range("A896").NumberFormat = "@"
$ws.Range("A896").Value = "2026/03/01"
$ws.Range("A896").Style = "Normal"
$ws.Range("B896").Value = "日"
$ws.Range("C896").Value = 1
$ws.Range("D896").Value = 201
